$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at row 909 (shifts existing rows 909:973 down to 912:976)
$ws.Range("A909:A911").EntireRow.Insert()

# Row 909 - Betarraga, Primera, week of 2022-07-04 (serial 44746)
$ws.Cells.Item(909, 1).Value = 6
$ws.Cells.Item(909, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(909, 3).Value = "Metropolitana"
$ws.Cells.Item(909, 4).Value = 44746
$ws.Cells.Item(909, 5).Value = 13
$ws.Cells.Item(909, 6).Value = 100114014
$ws.Cells.Item(909, 7).Value = "Betarraga"
$ws.Cells.Item(909, 8).Value = "Sin especificar"
$ws.Cells.Item(909, 9).Value = "Primera"
$ws.Cells.Item(909, 10).Value = 40000
$ws.Cells.Item(909, 11).Value = 110
$ws.Cells.Item(909, 12).Value = 120
$ws.Cells.Item(909, 13).Value = 115
$ws.Cells.Item(909, 14).Value = "`$/unidad"
$ws.Cells.Item(909, 15).Value = "Región Metropolitana"
$ws.Cells.Item(909, 16).Value = 115
$ws.Cells.Item(909, 17).Value = 1
$ws.Cells.Item(909, 18).Value = "Hortaliza"

# Row 910 - Betarraga, Segunda, same week
$ws.Cells.Item(910, 1).Value = 6
$ws.Cells.Item(910, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(910, 3).Value = "Metropolitana"
$ws.Cells.Item(910, 4).Value = 44746
$ws.Cells.Item(910, 5).Value = 13
$ws.Cells.Item(910, 6).Value = 100114014
$ws.Cells.Item(910, 7).Value = "Betarraga"
$ws.Cells.Item(910, 8).Value = "Sin especificar"
$ws.Cells.Item(910, 9).Value = "Segunda"
$ws.Cells.Item(910, 10).Value = 31000
$ws.Cells.Item(910, 11).Value = 85
$ws.Cells.Item(910, 12).Value = 95
$ws.Cells.Item(910, 13).Value = 90
$ws.Cells.Item(910, 14).Value = "`$/unidad"
$ws.Cells.Item(910, 15).Value = "Región Metropolitana"
$ws.Cells.Item(910, 16).Value = 90
$ws.Cells.Item(910, 17).Value = 1
$ws.Cells.Item(910, 18).Value = "Hortaliza"

# Row 911 - Betarraga, Tercera, same week
$ws.Cells.Item(911, 1).Value = 6
$ws.Cells.Item(911, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(911, 3).Value = "Metropolitana"
$ws.Cells.Item(911, 4).Value = 44746
$ws.Cells.Item(911, 5).Value = 13
$ws.Cells.Item(911, 6).Value = 100114014
$ws.Cells.Item(911, 7).Value = "Betarraga"
$ws.Cells.Item(911, 8).Value = "Sin especificar"
$ws.Cells.Item(911, 9).Value = "Tercera"
$ws.Cells.Item(911, 10).Value = 9000
$ws.Cells.Item(911, 11).Value = 70
$ws.Cells.Item(911, 12).Value = 70
$ws.Cells.Item(911, 13).Value = 70
$ws.Cells.Item(911, 14).Value = "`$/unidad"
$ws.Cells.Item(911, 15).Value = "Región Metropolitana"
$ws.Cells.Item(911, 16).Value = 70
$ws.Cells.Item(911, 17).Value = 1
$ws.Cells.Item(911, 18).Value = "Hortaliza"
